$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Helper: set a cell value as plain text, avoiding Excel's automatic
# numeric/date type coercion for strings that look like numbers
# (e.g. "30.686.16", "0.9972"), while leaving the cell style untouched.
function Set-TextValue($range, $text) {
    $range.NumberFormat = "@"
    $range.Value = $text
    $range.ClearFormats()
}

# Row 2: 'Bitcoin' -> 'Bitcoin'
Set-TextValue $ws.Range("D2") "30.686.16"
Set-TextValue $ws.Range("E2") "  +0.70%  "

# Row 3: 'Ethereum' -> 'Ethereum'
Set-TextValue $ws.Range("D3") "1.963.07"
Set-TextValue $ws.Range("E3") "  +2.63%  "

# Row 4: 'TetherUSD' -> 'TetherUSD'
Set-TextValue $ws.Range("D4") "0.9972"
Set-TextValue $ws.Range("E4") "  -0.15%  "

# Row 5: 'BNB' -> 'BNB'
Set-TextValue $ws.Range("D5") "248.62"
Set-TextValue $ws.Range("E5") "  +1.44%  "

# Row 6: 'USDC' -> 'USDC'
Set-TextValue $ws.Range("D6") "0.9979"
Set-TextValue $ws.Range("E6") "  -0.11%  "

# Row 7: 'XRP' -> 'XRP'
Set-TextValue $ws.Range("D7") "0.4809"
Set-TextValue $ws.Range("E7") "  -0.12%  "

# Row 8: 'Cardano' -> 'OKB'
$ws.Range("B8").Value = "OKB"
$ws.Range("C8").Value = "https://coinranking.com/coin/PDKcptVnzJTmN+okb-okb"
Set-TextValue $ws.Range("D8") "44.46"
Set-TextValue $ws.Range("E8") "  +0.12%  "

# Row 9: 'Dogecoin' -> 'Cardano'
$ws.Range("B9").Value = "Cardano"
$ws.Range("C9").Value = "https://coinranking.com/coin/qzawljRxB5bYu+cardano-ada"
Set-TextValue $ws.Range("D9") "0.2949"
Set-TextValue $ws.Range("E9") "  +1.99%  "

# Row 10: 'Litecoin' -> 'Dogecoin'
$ws.Range("B10").Value = "Dogecoin"
$ws.Range("C10").Value = "https://coinranking.com/coin/a91GCGd_u96cF+dogecoin-doge"
Set-TextValue $ws.Range("D10") "0.06807"
Set-TextValue $ws.Range("E10") "  +1.25%  "

# Row 11: 'Solana' -> 'Litecoin'
$ws.Range("B11").Value = "Litecoin"
$ws.Range("C11").Value = "https://coinranking.com/coin/D7B1x_ks7WhV5+litecoin-ltc"
Set-TextValue $ws.Range("D11") "112.10"
Set-TextValue $ws.Range("E11") "  +1.10%  "

# Row 12: 'WrappedEther' -> 'Solana'
$ws.Range("B12").Value = "Solana"
$ws.Range("C12").Value = "https://coinranking.com/coin/zNZHO_Sjf+solana-sol"
Set-TextValue $ws.Range("D12") "19.48"
Set-TextValue $ws.Range("E12") "  +2.30%  "

# Row 13: 'TRON' -> 'WrappedEther'
$ws.Range("B13").Value = "WrappedEther"
$ws.Range("C13").Value = "https://coinranking.com/coin/Mtfb0obXVh59u+wrappedether-weth"
Set-TextValue $ws.Range("D13") "1.956.74"
Set-TextValue $ws.Range("E13") "  +2.32%  "

# Row 14: 'Polkadot' -> 'TRON'
$ws.Range("B14").Value = "TRON"
$ws.Range("C14").Value = "https://coinranking.com/coin/qUhEFk1I61atv+tron-trx"
Set-TextValue $ws.Range("D14") "0.07693"
Set-TextValue $ws.Range("E14") "  +1.87%  "

# Row 15: 'Polygon' -> 'Polkadot'
$ws.Range("B15").Value = "Polkadot"
$ws.Range("C15").Value = "https://coinranking.com/coin/25W7FG7om+polkadot-dot"
Set-TextValue $ws.Range("D15") "5.491"
Set-TextValue $ws.Range("E15") "  +4.27%  "

# Row 16: 'BitcoinCash' -> 'Polygon'
$ws.Range("B16").Value = "Polygon"
$ws.Range("C16").Value = "https://coinranking.com/coin/uW2tk-ILY0ii+polygon-matic"
Set-TextValue $ws.Range("D16") "0.6877"
Set-TextValue $ws.Range("E16") "  +2.32%  "

# Row 17: 'WrappedBTC' -> 'BitcoinCash'
$ws.Range("B17").Value = "BitcoinCash"
$ws.Range("C17").Value = "https://coinranking.com/coin/ZlZpzOJo43mIo+bitcoincash-bch"
Set-TextValue $ws.Range("D17") "295.49"
Set-TextValue $ws.Range("E17") "  +2.73%  "

# Row 18: 'Avalanche' -> 'WrappedBTC'
$ws.Range("B18").Value = "WrappedBTC"
$ws.Range("C18").Value = "https://coinranking.com/coin/x4WXHge-vvFY+wrappedbtc-wbtc"
Set-TextValue $ws.Range("D18") "30.678.77"
Set-TextValue $ws.Range("E18") "  +0.64%  "

# Row 19: 'Uniswap' -> 'Avalanche'
$ws.Range("B19").Value = "Avalanche"
$ws.Range("C19").Value = "https://coinranking.com/coin/dvUj0CzDZ+avalanche-avax"
Set-TextValue $ws.Range("D19") "13.26"
Set-TextValue $ws.Range("E19") "  +3.21%  "

# Row 20: 'WrappedliquidstakedEther2.0' -> 'Uniswap'
$ws.Range("B20").Value = "Uniswap"
$ws.Range("C20").Value = "https://coinranking.com/coin/_H5FVG9iW+uniswap-uni"
Set-TextValue $ws.Range("D20") "5.668"
Set-TextValue $ws.Range("E20") "  +3.58%  "

# Row 21: 'ShibaInu' -> 'WrappedliquidstakedEther2.0'
$ws.Range("B21").Value = "WrappedliquidstakedEther2.0"
$ws.Range("C21").Value = "https://coinranking.com/coin/CiixT63n3+wrappedliquidstakedether20-wsteth"
Set-TextValue $ws.Range("D21") "2.229.24"
Set-TextValue $ws.Range("E21") "  +3.05%  "

# Row 22: 'Dai' -> 'ShibaInu'
$ws.Range("B22").Value = "ShibaInu"
$ws.Range("C22").Value = "https://coinranking.com/coin/xz24e0BjL+shibainu-shib"
Set-TextValue $ws.Range("D22") "0.000007679"
Set-TextValue $ws.Range("E22") "  +1.07%  "

# Row 23: 'BinanceUSD' -> 'Dai'
$ws.Range("B23").Value = "Dai"
$ws.Range("C23").Value = "https://coinranking.com/coin/MoTuySvg7+dai-dai"
Set-TextValue $ws.Range("D23") "0.9985"
Set-TextValue $ws.Range("E23") "  -0.04%  "

# Row 24: 'Chainlink' -> 'BinanceUSD'
$ws.Range("B24").Value = "BinanceUSD"
$ws.Range("C24").Value = "https://coinranking.com/coin/vSo2fu9iE1s0Y+binanceusd-busd"
Set-TextValue $ws.Range("D24") "0.9947"
Set-TextValue $ws.Range("E24") "  -0.38%  "

# Row 25: 'Cosmos' -> 'Chainlink'
$ws.Range("B25").Value = "Chainlink"
$ws.Range("C25").Value = "https://coinranking.com/coin/VLqpJwogdhHNb+chainlink-link"
Set-TextValue $ws.Range("D25") "6.609"
Set-TextValue $ws.Range("E25") "  +3.10%  "

# Row 26: 'Monero' -> 'Cosmos'
$ws.Range("B26").Value = "Cosmos"
$ws.Range("C26").Value = "https://coinranking.com/coin/Knsels4_Ol-Ny+cosmos-atom"
Set-TextValue $ws.Range("D26") "9.751"
Set-TextValue $ws.Range("E26") "  +3.04%  "

# Row 27: 'EthereumClassic' -> 'Monero'
$ws.Range("B27").Value = "Monero"
$ws.Range("C27").Value = "https://coinranking.com/coin/3mVx2FX_iJFp5+monero-xmr"
Set-TextValue $ws.Range("D27") "169.01"
Set-TextValue $ws.Range("E27") "  +3.11%  "

# Row 28: 'LidoDAOToken' -> 'EthereumClassic'
$ws.Range("B28").Value = "EthereumClassic"
$ws.Range("C28").Value = "https://coinranking.com/coin/hnfQfsYfeIGUQ+ethereumclassic-etc"
Set-TextValue $ws.Range("D28") "20.32"
Set-TextValue $ws.Range("E28") "  +0.16%  "

# Row 29: 'Stellar' -> 'LidoDAOToken'
$ws.Range("B29").Value = "LidoDAOToken"
$ws.Range("C29").Value = "https://coinranking.com/coin/Pe93bIOD2+lidodaotoken-ldo"
Set-TextValue $ws.Range("D29") "2.211"
Set-TextValue $ws.Range("E29") "  +4.63%  "

# Row 30: 'Toncoin' -> 'Stellar'
$ws.Range("B30").Value = "Stellar"
$ws.Range("C30").Value = "https://coinranking.com/coin/f3iaFeCKEmkaZ+stellar-xlm"
Set-TextValue $ws.Range("D30") "0.1089"
Set-TextValue $ws.Range("E30") "  +3.31%  "

# Row 31: 'Filecoin' -> 'Toncoin'
$ws.Range("B31").Value = "Toncoin"
$ws.Range("C31").Value = "https://coinranking.com/coin/67YlI0K1b+toncoin-ton"
Set-TextValue $ws.Range("D31") "1.441"
Set-TextValue $ws.Range("E31") "  +2.66%  "

# Row 32: 'InternetComputer(DFINITY)' -> 'Filecoin'
$ws.Range("B32").Value = "Filecoin"
$ws.Range("C32").Value = "https://coinranking.com/coin/ymQub4fuB+filecoin-fil"
Set-TextValue $ws.Range("D32") "4.660"
Set-TextValue $ws.Range("E32") "  +15.22%  "

# Row 33: 'Hedera' -> 'InternetComputer(DFINITY)'
$ws.Range("B33").Value = "InternetComputer(DFINITY)"
$ws.Range("C33").Value = "https://coinranking.com/coin/aMNLwaUbY+internetcomputerdfinity-icp"
Set-TextValue $ws.Range("D33") "4.440"
Set-TextValue $ws.Range("E33") "  +6.39%  "

# Row 34: 'ImmutableX' -> 'Hedera'
$ws.Range("B34").Value = "Hedera"
$ws.Range("C34").Value = "https://coinranking.com/coin/jad286TjB+hedera-hbar"
Set-TextValue $ws.Range("D34") "0.05081"
Set-TextValue $ws.Range("E34") "  +2.08%  "

# Row 35: 'ARBITRUM' -> 'ImmutableX'
$ws.Range("B35").Value = "ImmutableX"
$ws.Range("C35").Value = "https://coinranking.com/coin/Z96jIvLU7+immutablex-imx"
Set-TextValue $ws.Range("D35") "0.7759"
Set-TextValue $ws.Range("E35") "  +6.53%  "

# Row 36: 'VeChain' -> 'ARBITRUM'
$ws.Range("B36").Value = "ARBITRUM"
$ws.Range("C36").Value = "https://coinranking.com/coin/1Uo6s62Oc+arbitrum-arb"
Set-TextValue $ws.Range("D36") "1.172"
Set-TextValue $ws.Range("E36") "  +3.40%  "

# Row 37: 'HuobiToken' -> 'VeChain'
$ws.Range("B37").Value = "VeChain"
$ws.Range("C37").Value = "https://coinranking.com/coin/FEbS54wxo4oIl+vechain-vet"
Set-TextValue $ws.Range("D37") "0.02074"
Set-TextValue $ws.Range("E37") "  +2.01%  "

# Row 38: 'MXToken' -> 'HuobiToken'
$ws.Range("B38").Value = "HuobiToken"
$ws.Range("C38").Value = "https://coinranking.com/coin/DXwP4wF9ksbBO+huobitoken-ht"
Set-TextValue $ws.Range("D38") "2.729"
Set-TextValue $ws.Range("E38") "  +0.25%  "

# Row 39: 'RenderToken' -> 'MXToken'
$ws.Range("B39").Value = "MXToken"
$ws.Range("C39").Value = "https://coinranking.com/coin/QUC5kVAxSoB-+mxtoken-mx"
Set-TextValue $ws.Range("D39") "2.710"
Set-TextValue $ws.Range("E39") "  +1.69%  "

# Row 40: 'Quant' -> 'RenderToken'
$ws.Range("B40").Value = "RenderToken"
$ws.Range("C40").Value = "https://coinranking.com/coin/7C4Mh4xy1yDel+rendertoken-rndr"
Set-TextValue $ws.Range("D40") "2.072"
Set-TextValue $ws.Range("E40") "  +3.02%  "

# Row 41: 'TheSandbox' -> 'Quant'
$ws.Range("B41").Value = "Quant"
$ws.Range("C41").Value = "https://coinranking.com/coin/bauj_21eYVwso+quant-qnt"
Set-TextValue $ws.Range("D41") "111.07"
Set-TextValue $ws.Range("E41") "  +0.46%  "

# Row 42: 'TrustWalletToken' -> 'TheSandbox'
$ws.Range("B42").Value = "TheSandbox"
$ws.Range("C42").Value = "https://coinranking.com/coin/pxtKbG5rg+thesandbox-sand"
Set-TextValue $ws.Range("D42") "0.4469"
Set-TextValue $ws.Range("E42") "  +0.85%  "

# Row 43: 'FraxShare' -> 'TrustWalletToken'
$ws.Range("B43").Value = "TrustWalletToken"
$ws.Range("C43").Value = "https://coinranking.com/coin/Hm3OlynlC+trustwallettoken-twt"
Set-TextValue $ws.Range("D43") "0.8734"
Set-TextValue $ws.Range("E43") "  +0.70%  "

# Row 44: 'PaxDollar' -> 'FraxShare'
$ws.Range("B44").Value = "FraxShare"
$ws.Range("C44").Value = "https://coinranking.com/coin/3nNpuxHJ8+fraxshare-fxs"
Set-TextValue $ws.Range("D44") "6.004"
Set-TextValue $ws.Range("E44") "  +3.15%  "

# Row 45: 'Aave' -> 'Aave'
Set-TextValue $ws.Range("D45") "69.92"
Set-TextValue $ws.Range("E45") "  +2.47%  "

# Row 46: 'Aptos' -> 'PaxDollar'
$ws.Range("B46").Value = "PaxDollar"
$ws.Range("C46").Value = "https://coinranking.com/coin/JCKLgWPAF+paxdollar-usdp"
Set-TextValue $ws.Range("D46") "0.9992"
Set-TextValue $ws.Range("E46") "  +0.02%  "

# Row 47: 'EnergySwap' -> 'Aptos'
$ws.Range("B47").Value = "Aptos"
$ws.Range("C47").Value = "https://coinranking.com/coin/HGYj5JCv5+aptos-apt"
Set-TextValue $ws.Range("D47") "7.411"
Set-TextValue $ws.Range("E47") "  +1.05%  "

# Row 48: 'Algorand' -> 'EnergySwap'
$ws.Range("B48").Value = "EnergySwap"
$ws.Range("C48").Value = "https://coinranking.com/coin/SbWqqTui-+energyswap-ens"
Set-TextValue $ws.Range("D48") "9.391"
Set-TextValue $ws.Range("E48") "  +0.90%  "

# Row 49: 'BitcoinSV' -> 'Algorand'
$ws.Range("B49").Value = "Algorand"
$ws.Range("C49").Value = "https://coinranking.com/coin/TpHE2IShQw-sJ+algorand-algo"
Set-TextValue $ws.Range("D49") "0.1254"
Set-TextValue $ws.Range("E49") "  +1.17%  "

# Row 50: 'Elrond' -> 'BitcoinSV'
$ws.Range("B50").Value = "BitcoinSV"
$ws.Range("C50").Value = "https://coinranking.com/coin/VcMY11NONHSA0+bitcoinsv-bsv"
Set-TextValue $ws.Range("D50") "48.07"
Set-TextValue $ws.Range("E50") "  -1.85%  "

# Row 51: 'WOONetwork' -> 'Elrond'
$ws.Range("B51").Value = "Elrond"
$ws.Range("C51").Value = "https://coinranking.com/coin/omwkOTglq+elrond-egld"
Set-TextValue $ws.Range("D51") "35.70"
Set-TextValue $ws.Range("E51") "  +2.64%  "
